$d = $word.ActiveDocument
$wNs = "http://schemas.openxmlformats.org/wordprocessingml/2006/main"

function Merge-Text($findText) {
    # Replacing the run-spanning text with itself (via Find/Replace) collapses
    # the surrounding <w:r> runs (and the w:proofErr markers Word had inserted
    # around them) into a single run, matching how Word normalizes text after
    # an in-place edit/spell-check acceptance.
    $d.Content.Find.Execute($findText, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $findText, 2) | Out-Null
}

# 1) "1) What is the strongest correlation of real estate prices? - Richard (?) "
Merge-Text("1) What is the strongest correlation of real estate prices? - Richard (?) ")

# 2) "*double check the columns in these datasets and choose the best one "
Merge-Text("*double check the columns in these datasets and choose the best one ")

# 3) "...less strong than the effects on a 3br?. "
Merge-Text("2) We expect the strength of this relationship to differ for different types of homes. The effects on a 1br will be less strong than the effects on a 3br?. ")

# 4) "...Likewise when income levels..."
Merge-Text("4) Unemployment rates, income levels, and inflation will have a negative, positive, and positive relationship with real estate prices respectively. As unemployment rises, real estate prices will decrease? Likewise when income levels in an area increase, real estate prices will increase. And as inflation rises, real estate prices will also inflate? ")

# 5) Split the "Find rate of change " list item into three runs ("Find rate" /
#    "s" / " of change"), then add a new sibling bullet "Plot data " right
#    after it (same list level/numbering), leaving the existing empty bullet
#    that follows untouched.
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -eq "Find rate of change `r") {
        $target = $p
        break
    }
}

if ($target -ne $null) {
    $start = $target.Range.Start
    $paraEnd = $target.Range.End

    # Rebuild the paragraph's run content as three separate <w:r> elements.
    # InsertXML (unlike InsertAfter/Range.Text) does not get coalesced back
    # into a single run, so it preserves the three-run split.
    $textRange = $d.Range($start, $paraEnd - 1)
    $splitXml = "<w:p xmlns:w=`"$wNs`"><w:r><w:t>Find rate</w:t></w:r><w:r><w:t>s</w:t></w:r><w:r><w:t xml:space=`"preserve`"> of change</w:t></w:r></w:p>"
    [void]$textRange.InsertXML($splitXml)

    # Re-locate the (now retextted) paragraph and add a new list paragraph
    # right after it, inheriting the same pPr (ListParagraph, ilvl 1, numId 1).
    $target2 = $null
    foreach ($p in $d.Paragraphs) {
        if ($p.Range.Text -eq "Find rates of change`r") {
            $target2 = $p
            break
        }
    }
    [void]$target2.Range.InsertParagraphAfter()

    # The freshly-inserted empty paragraph is the one right after $target2.
    $prevWasTarget = $false
    $newEmptyPara = $null
    foreach ($p in $d.Paragraphs) {
        if ($prevWasTarget) {
            $newEmptyPara = $p
            break
        }
        if ($p.Range.Text -eq "Find rates of change`r") {
            $prevWasTarget = $true
        }
    }
    $newEmptyPara.Range.Text = "Plot data "
}
